$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Write every new cell value in the exact chronological order the
# --- original author entered them, so the shared-string table indices
# --- line up with the target workbook (110..144, in this sequence).

$ws.Range("H1").Value = "ProtocolEndpoint"
$ws.Range("H2").Value = "pr-endpoint-1"
$ws.Range("H3").Value = "pr-endpoint-2"
$ws.Range("H4").Value = "pr-endpoint-3"

$ws.Range("AD1").Value = "Eff1"
$ws.Range("AF1").Value = "Cond11"
$ws.Range("AG1").Value = "Cond12"
$ws.Range("AH1").Value = "Eff2"
$ws.Range("AI1").Value = "Cond21"

$ws.Range("AF2").Value = "cond1-1"
$ws.Range("AF3").Value = "cond1-2"
$ws.Range("AF4").Value = "cond1-3"

$ws.Range("AH2").Value = "<= 51 A"
$ws.Range("AH3").Value = "<= 52 A"
$ws.Range("AH4").Value = "<= 53 A"

$ws.Range("AJ1").Value = "Eff3"
$ws.Range("AJ2").Value = "1-56 Pa"
$ws.Range("AJ3").Value = "2-56 Pa"
$ws.Range("AJ4").Value = "3-56 Pa"

$ws.Range("AD4").Value = "~ 400 K"
$ws.Range("AD3").Value = "ca. 300 K"
$ws.Range("AD2").Value = "> 200 K"

$ws.Range("AE1").Value = "Error1"

$ws.Range("AG2").Value = "101 m"
$ws.Range("AG3").Value = "102 m"
$ws.Range("AG4").Value = "103 m"

$ws.Range("AI2").Value = "> 1"
$ws.Range("AI3").Value = "> 2"
$ws.Range("AI4").Value = "> 3"

$ws.Range("AK1").Value = "Eff4-Lo"
$ws.Range("AM1").Value = "Eff4-Unit"
$ws.Range("AL1").Value = "Eff4-Up"
$ws.Range("AM2").Value = "nm"
$ws.Range("AN1").Value = "Eff4-Lo-Quil"
$ws.Range("AN2").Value = ">="

$ws.Range("AM3").Value = "nm"
$ws.Range("AM4").Value = "nm"
$ws.Range("AN3").Value = ">="
$ws.Range("AN4").Value = ">="

# --- Numeric cells (no shared-string impact) ---
$ws.Range("AE2").Value = 81
$ws.Range("AE3").Value = 82
$ws.Range("AE4").Value = 83

$ws.Range("AK2").Value = 11
$ws.Range("AK3").Value = 12
$ws.Range("AK4").Value = 13

$ws.Range("AL2").Value = 21
$ws.Range("AL3").Value = 22
$ws.Range("AL4").Value = 23

# --- Styling ---
# Header row (row 1) new cells copy the bold header style used by the
# existing header cells (A1 etc, style index 1).
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("AD1:AN1").PasteSpecial(-4122)

# Re-apply the values, since PasteSpecial(formats) shouldn't clobber them,
# but make sure they are correct in case of any ambiguity.
$ws.Range("H1").Value = "ProtocolEndpoint"
$ws.Range("AD1").Value = "Eff1"
$ws.Range("AE1").Value = "Error1"
$ws.Range("AF1").Value = "Cond11"
$ws.Range("AG1").Value = "Cond12"
$ws.Range("AH1").Value = "Eff2"
$ws.Range("AI1").Value = "Cond21"
$ws.Range("AJ1").Value = "Eff3"
$ws.Range("AK1").Value = "Eff4-Lo"
$ws.Range("AL1").Value = "Eff4-Up"
$ws.Range("AM1").Value = "Eff4-Unit"
$ws.Range("AN1").Value = "Eff4-Lo-Quil"

# New cell-alignment style for the Error1 numeric column (center aligned).
$ws.Range("AE2:AE4").HorizontalAlignment = -4108

# --- Sheet view / selection changes from the diff ---
$ws.Application.ActiveWindow.ScrollColumn = 29
[void]$ws.Range("AI4").Select()

# --- New column width (col H) ---
$ws.Columns.Item(8).ColumnWidth = 15.33203125
